$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "Asking for verbal consent is similar to requesting ID or issuing a receipt, which involves direct interaction with the customer. Considering the need for clear communication and potential customer queries, this activity might take slightly longer than issuing a receipt but less than reviewing an ID. Given the technical infrastructure limitations and the need for staff to ensure the customer fully understands the consent process, this could realistically take around 15 seconds."

$ws.Range("C5").Value = "Requesting a FIN ID from the applicant is similar to requesting an ID, which takes approximately 14.55 seconds. Considering the technical infrastructure in Ethiopia, staff training, and potential customer interaction time, a slight increase to 15 seconds accounts for any minor delays or retries that might occur during the process."

$ws.Range("C6").Value = "This activity would take approximately 15 seconds because scanning a QR code is generally a quick process similar to scanning a barcode, which takes about 23.4 seconds in the current process. However, QR codes are typically faster to scan than barcodes. Considering the technical infrastructure limitations in Ethiopia, potential staff training needs, and customer interaction time, a slightly reduced time is estimated. Additionally, potential delays or retries are accounted for, but are expected to be minimal due to the straightforward nature of QR code scanning."

$ws.Range("C7").Value = "This activity would take approximately 45 seconds because it is similar to the `"Fill data into EKYC System`" activity, which takes 53.59 seconds. Considering the technical infrastructure limitations in Ethiopia, staff training needs, and potential customer interaction time, the process might be slightly faster due to familiarity with the system but still requires significant time for accuracy and potential retries."

$ws.Range("C9").Value = "Linking the unique barcode on the SIM card with the FIN is a straightforward technical task similar to scanning the barcode on the SIM card, which takes about 23.4 seconds. However, this new activity might be slightly quicker as it involves linking rather than scanning. Considering the technical infrastructure limitations in Ethiopia, potential staff training needs, and customer interaction time, this task could be completed in approximately 15 seconds. This estimate accounts for potential delays or retries due to system lags or errors."

$ws.Range("D9").Value = 15
